# Generate Report for Handback
#
# The "f95e20a6-b991-4ddd-b81f-a46dffff057c" file has now been handed back
# (it was previously only "Ready for handoff"). Update its Status on every
# sheet, and record the new "Latest Handback DateTime" on the per-language
# sheets.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- zh-cn sheet: file f95e20a6... is row 3 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("H3").Value = "2016-03-12 02:33:09"

# --- de-de sheet: file f95e20a6... is row 3 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("H3").Value = "2016-03-12 02:33:14"

# --- Overview sheet: reflect the updated status for both locales ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus
